$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$titleShape = $s.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange

$newText = "9 Assignment Operator"

# Rewrite the title's first run in place (this preserves that run's own
# rPr - lang/altLang/sz/b/dirty/solidFill - exactly as authored) so the
# merged run ends up with the same formatting as the original "9 " run.
$tr.Runs(1).Text = $newText

# The first run now contains the full desired text, followed by the
# leftover text of the runs that used to hold "Assignmnent" and
# " Operator". Select that leftover as a character range and clear it so
# the now-redundant run elements are removed entirely (not just emptied).
$fullLen = $tr.Text.Length
$newLen = $newText.Length
if ($fullLen -gt $newLen) {
    $leftover = $tr.Characters($newLen + 1, $fullLen - $newLen)
    $leftover.Text = ""
}
